# article 20 is live
# Adds row 21 (article/Ser 20) to Sheet1 with date, ayats, content, author and tags,
# matching the structure of the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ser (article number)
$ws.Cells.Item(21, 1).Value = 20

# Date (2020-01-20), formatted like the date above it
$ws.Cells.Item(21, 2).Value = 43850
$ws.Cells.Item(21, 2).NumberFormat = "d-mmm-yy"

# Ayats
$ws.Cells.Item(21, 3).Value = "Surah Baqarah, 211 - 217"

# Tags
$ws.Cells.Item(21, 6).Value = "Battle of trench, Spending wealth, help of Allah (swt)"

# Content
$ws.Cells.Item(21, 4).Value = "h3: Battle of trench.`np: As people migrated to Madina from Makkah, they left their wives, kids, families, trades behind. Poverty on their heads, heads facing down, sun scorching above, desert eating below, feet trembling they moved to find Allah in the emptiness. Stakes were high. Food they do not find, they die. Shelter they do not find, they die. Yet they moved with all risks of losing. They had lost most of their lives. Everything was at stake. One Muhammad (saw) is beside them, they have not seen Allah (swt), did not see angels in real. `np: What moved them during migration? What kept them moving through 12 years of sufferings? Moving is easy if it is a high roof bus with food served on table AC working fine. But when your friends fall beside you, heat eats your brains, death is around the corner, moving is tough. `nh3: What moved them..`np: The trust in Allah (swt) and His messenger (PBUH). `np: Did hardships end in their lives? No. After this phase another phase started. `np: In today’s lesson Allah (swt) instructed muslims to stay steadfast during the battle of trench. `nquote: Or do you think that you will enter Paradise while such [trial] has not yet come to you as came to those who passed on before you? They were touched by poverty and hardship and were shaken until [even their] messenger and those who believed with him said,`"When is the help of Allah ?`" Unquestionably, the help of Allah is near. <br> - Surah Baqarah verse 214`np: The help of Allah is near for people who are facing calamities.`nh3: Simple message.`np: <b>1. Beautified this world (ayat 211) </b>for people who deviates from truth. `np: <b>2. Spend wealth (ayat 215)</b> on parents, relatives, orphans, needy and the travellers.`np: <b>3. Fight the fitna (ayat 217)</b> as it is greater then killing. `np: <b>4. Do not revert back (ayat 217)</b> after coming to the path of Allah (swt). `nquote: And whoever of you reverts from his religion [to disbelief] and dies while he is a disbeliever - for those, their deeds have become worthless in this world and the Hereafter, and those are the companions of the Fire, they will abide therein eternally. <br> - Surah Baqarh verse 217`np: Poetry or no poetry, emotional or not, people will follow Allah (swt) with all their energy. People who do not want to follow, they will not. Quran is just a book. Messengers are just people. Satan is just a jin. We want to get lost in fitna? We will get lost. Evading death is not possible. Running before death to find the door to freedom is useless. Death is a shackle we can not break, yet we are running for freedom. <b class='lavendar'>Be free</b>, in the kingdom of Allah (swt). `np.note: Guys, if you like this project. Please follow this project's page on twitter. <a href=`"https://twitter.com/zakatlists`">Click here to go to the twitter page</a>."

# Author
$ws.Cells.Item(21, 5).Value = "Qasim Ali"

# Match row height of the other fully-populated rows
$ws.Rows.Item(21).RowHeight = 409.6

# Scroll / select the newly-added row like the author left it
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D21").Select()
